# Splits two single-run paragraphs ("Norma de recuperacao" note and the
# "Bibliografia" entry) into multiple <w:t> segments separated by manual
# line breaks (<w:br/>), one per sentence / reference entry, while keeping
# everything inside a single run (so existing run formatting is preserved).
#
# Word's Find/Replace treats "^l" in the replacement text as a manual line
# break (w:br), and when MatchWildcards is $false the Find What text is
# matched literally (parentheses, '+', etc. are not treated as regex).

$d = $word.ActiveDocument

$find1 = 'MF = (0,5 M + 0,5 R)M = Média de aproveitamento do aluno, antes da recuperaçãoR = Nota de uma prova de recuperaçãoMF = nota final de aproveitamento, após a recuperaçãoAprovação com média final de aproveitamento maior ou igual a 5,0.A recuperação deverá consistir de uma prova escrita englobando a matéria toda do semestre.Terá direito à prova de recuperação aqueles alunos reprovados com nota acima de 3,0 e frequência mínima de 70%.'
$repl1 = 'MF = (0,5 M + 0,5 R)^lM = Média de aproveitamento do aluno, antes da recuperação^lR = Nota de uma prova de recuperação^lMF = nota final de aproveitamento, após a recuperação^lAprovação com média final de aproveitamento maior ou igual a 5,0.^lA recuperação deverá consistir de uma prova escrita englobando a matéria toda do semestre.^lTerá direito à prova de recuperação aqueles alunos reprovados com nota acima de 3,0 e frequência mínima de 70%.'

$find2 = 'BARNEY, J.B.; CLARK, D.N. Resource-Based Theory: Creating and Sustaining Competitive Advantage. Oxford University Press, 2007.BESSANT, J.; TIDD, J. Inovação e empreendedorismo. Porto Alegre, Bookman, 2009.BURGELMAN, R. A.; CHRISTENSEN, C. M.; WHEELWRIGTH, S. C. Gestão estratégica da tecnologia e da inovação: conceitos e soluções. AMGH Editora, 2013.CONWAY, S; STEWARD, F. Managing and shaping innovation. Oxford University Press, 2009.CHRISTENSEN, Clayton M. O dilema da inovação. São Paulo: Makron Books, 2011.DAVILA, T; EPSTEIN, M. J.; SHELTON, R. As regras da Inovação. Porto Alegre, Bookman, 2008.DE NEGRI, J.A; SALERNO, M.S., (Orgs.). Inovação, padrões tecnológicos e desempenho das firmas industriais brasileiras. Brasília, Ipea, 2005.DODGSON, M.; GANN, D.; SALTER, A. The management of technological innovation: strategy and practice. Oxford University Press, 2008.DRUCKER, P.F. Inovação e espírito empreendedor. São Paulo: Pioneira, 2000.FIGUEIREDO, P.N. Gestão da inovação: conceitos, métricas e experiências de empresas no Brasil. Rio de Janeiro, LTC, 2015.FITZGERALD, E. et al. Inside Real Innovation: How the Right Approach Can Move Ideas from R&D to Market-And Get the Economy Moving. World Scientific, 2011.GOFFIN, K.; MITCHELL, R. Innovation management. 3nd ed. Palgrave – MacMillan, Houndsmill, 2017.HELFAT, C.E. et al. Dynamic capabilities: understanding strategic change in organizations. Blackwell Publishing, 2007.PRAHALAD,C.K.; KRISHNAN,M.S. The new of innovation. EUA: Editora Soundview Executive Book Sumaries, 2008.PROENÇA, A. et al. Gestão da inovação e competitividade no Brasil: da teoria para a prática. Bookman Editora, 2015.SALERNO, M.S.; GOMES, L.A.V. Gestão da inovação (mais) radical. Rio de Janeiro: Elsevier, 2018.SCHILLING, M.A. Strategic management of technological innovation. MacGraw-Hill/Irwin, 2009. TEECE, D. Capabilities and strategic management. In: Edited by Foss. N. Resources firms and strategies. A reader in the Resource-based Perspective. WA: Ed. Oxford University, 1987.TIDD, J.; BESSANT, J. Gestão da Inovação. Porto Alegre, Bookman, 2015.TIDD, J.; BESSANT, J. Strategic innovation management, Wiley, 2014.TIGRE, P. B. Gestão da inovação. Rio de Janeiro, Campus-Elsevier, 2006.TROTT, P. innovation management and new product development. Prentice Hall, 2008.WHITE, M. A.; BRUTON, G.D. The management of technology and innovation: a strategic approach. South-Western, Cengage Learning, 2011.'
$repl2 = 'BARNEY, J.B.; CLARK, D.N. Resource-Based Theory: Creating and Sustaining Competitive Advantage. Oxford University Press, 2007.^lBESSANT, J.; TIDD, J. Inovação e empreendedorismo. Porto Alegre, Bookman, 2009.^lBURGELMAN, R. A.; CHRISTENSEN, C. M.; WHEELWRIGTH, S. C. Gestão estratégica da tecnologia e da inovação: conceitos e soluções. AMGH Editora, 2013.^lCONWAY, S; STEWARD, F. Managing and shaping innovation. Oxford University Press, 2009.^lCHRISTENSEN, Clayton M. O dilema da inovação. São Paulo: Makron Books, 2011.^lDAVILA, T; EPSTEIN, M. J.; SHELTON, R. As regras da Inovação. Porto Alegre, Bookman, 2008.^lDE NEGRI, J.A; SALERNO, M.S., (Orgs.). Inovação, padrões tecnológicos e desempenho das firmas industriais brasileiras. Brasília, Ipea, 2005.^lDODGSON, M.; GANN, D.; SALTER, A. The management of technological innovation: strategy and practice. Oxford University Press, 2008.^lDRUCKER, P.F. Inovação e espírito empreendedor. São Paulo: Pioneira, 2000.^lFIGUEIREDO, P.N. Gestão da inovação: conceitos, métricas e experiências de empresas no Brasil. Rio de Janeiro, LTC, 2015.^lFITZGERALD, E. et al. Inside Real Innovation: How the Right Approach Can Move Ideas from R&D to Market-And Get the Economy Moving. World Scientific, 2011.^lGOFFIN, K.; MITCHELL, R. Innovation management. 3nd ed. Palgrave – MacMillan, Houndsmill, 2017.^lHELFAT, C.E. et al. Dynamic capabilities: understanding strategic change in organizations. Blackwell Publishing, 2007.^lPRAHALAD,C.K.; KRISHNAN,M.S. The new of innovation. EUA: Editora Soundview Executive Book Sumaries, 2008.^lPROENÇA, A. et al. Gestão da inovação e competitividade no Brasil: da teoria para a prática. Bookman Editora, 2015.^lSALERNO, M.S.; GOMES, L.A.V. Gestão da inovação (mais) radical. Rio de Janeiro: Elsevier, 2018.^lSCHILLING, M.A. Strategic management of technological innovation. MacGraw-Hill/Irwin, 2009. ^lTEECE, D. Capabilities and strategic management. In: Edited by Foss. N. Resources firms and strategies. A reader in the Resource-based Perspective. WA: Ed. Oxford University, 1987.^lTIDD, J.; BESSANT, J. Gestão da Inovação. Porto Alegre, Bookman, 2015.^lTIDD, J.; BESSANT, J. Strategic innovation management, Wiley, 2014.^lTIGRE, P. B. Gestão da inovação. Rio de Janeiro, Campus-Elsevier, 2006.^lTROTT, P. innovation management and new product development. Prentice Hall, 2008.^lWHITE, M. A.; BRUTON, G.D. The management of technology and innovation: a strategic approach. South-Western, Cengage Learning, 2011.'

function Replace-Literal($findText, $replText) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Replacement.ClearFormatting()
    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
    #          MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)
    return $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replText, 2)
}

$ok1 = Replace-Literal $find1 $repl1
Write-Output "Norma de recuperacao split: $ok1"
if (-not $ok1) { throw "Could not find the 'Norma de recuperacao' run to split." }

$ok2 = Replace-Literal $find2 $repl2
Write-Output "Bibliografia split: $ok2"
if (-not $ok2) { throw "Could not find the 'Bibliografia' run to split." }
